$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39 (entry "C-1") is superseded by a newly published entry ---
# Old: C-1 / "Climate change poses underappreciated threat to mesic
#      forests" / status "Revise and resubmit" / prev_status 2
# New: B-38 / "Drought sensitivity in mesic forests heightens their
#      vulnerability to climate change" / status "Peer-reviewed article" /
#      prev_status 3
$ws.Range("A39").Value = "B-38"
$ws.Range("C39").Value = "Drought sensitivity in mesic forests heightens their vulnerability to climate change"
$ws.Range("H39").Value = "Peer-reviewed article"
$ws.Range("K39").Value = 3

# --- Rows 40-43 keep their data but renumber the "C-n" series labels ---
# since the old C-1 entry was removed from that numbering sequence:
# C-2 -> C-1, C-3 -> C-2, C-4 -> C-3, C-5 -> C-4
$ws.Range("A40").Value = "C-1"
$ws.Range("A41").Value = "C-2"
$ws.Range("A42").Value = "C-3"
$ws.Range("A43").Value = "C-4"

# --- Refresh the sheet's viewport/selection to the edited area ---
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C40").Select()
